# Fill columns G:N with "N/A" for the rows that currently only have data
# through column F. This brings those rows in line with the "full" rows
# that already have values in every column A:N.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(5, 6, 10, 14, 15, 18, 22, 23, 27, 32, 33, 34, 37, 38, 40, 45, 50, 54, 55, 56, 61, 66, 67, 71, 72, 79, 80, 82, 83, 88, 93, 96, 103, 106, 112, 113, 114, 115, 119, 123, 125, 130, 131, 137, 138, 141, 142, 144, 145, 149, 156, 159, 161, 163, 165, 167, 172, 178, 181, 183, 184, 188, 189, 191)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
